# Generate Report for Handoff
# Updates the localization-status report: moves the "ee778694-...md" item
# from "In Translation" to "Ready for handoff" and refreshes the HO Xliff
# generate / handoff timestamps across the Overview, zh-cn and de-de
# sheets. Also widens the "Status"-adjacent date columns to fit the new
# values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2/F2 = zh-cn / de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 20:57:28"

# --- zh-cn sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 20:57:24"

# --- de-de sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 20:57:28"

# --- Column widths -------------------------------------------------------
# The zh-cn / de-de status columns (Overview!E:F, zh-cn!C, de-de!C) grow
# to fit the longer "Ready for handoff" label.
$newColumnWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
